{"js": "// Paragraph 4 currently ends with:\n//   \"...bazuj\u0105c na pocz\u0105tku i ko\u0144cu sekwencji.\"\n// We need it to end with:\n//   \"...bazuj\u0105c na pocz\u0105tku i ko\u0144cu sekwencji oraz chromosomie.\"\n// split across two runs (same rPr/lang=\"pl-PL\"):\n//   run 1: \"...sekwencji\"\n//   run 2: \" oraz chromosomie.\"   (xml:space=\"preserve\")\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph that contains the target sentence.\nconst needle = \"bazuj\u0105c na pocz\u0105tku i ko\u0144cu sekwencji.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(needle) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the target paragraph.\");\n}\n\n// Step 1: replace \"sekwencji.\" with \"sekwencji oraz chromosomie.\" as a single\n// edit so the run keeps uniform formatting (this naturally collapses into\n// one merged run).\nconst periodMatches = target.search(\"sekwencji.\", { matchCase: true });\nawait context.sync();\nif (periodMatches.items.length === 0) {\n  throw new Error(\"Could not find 'sekwencji.' to replace.\");\n}\nperiodMatches.items[0].insertText(\"sekwencji oraz chromosomie.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Step 2: re-select just the newly appended suffix (\" oraz chromosomie.\")\n// and nudge a character-formatting property on/off. This forces the host to\n// materialize it as its own run (matching the target OOXML: two runs with\n// identical rPr) instead of silently re-merging it with the previous run.\nconst suffixMatches = target.search(\" oraz chromosomie.\", { matchCase: true });\nawait context.sync();\nif (suffixMatches.items.length === 0) {\n  throw new Error(\"Could not find the newly inserted suffix.\");\n}\nconst suffixRange = suffixMatches.items[0];\nsuffixRange.font.bold = true;\nawait context.sync();\nsuffixRange.font.bold = false;\nawait context.sync();\n", "ps1": "# Paragraph 4 currently ends with:\n#   \"...bazuj\u0105c na pocz\u0105tku i ko\u0144cu sekwencji.\"\n# We need it to end with:\n#   \"...bazuj\u0105c na pocz\u0105tku i ko\u0144cu sekwencji oraz chromosomie.\"\n# split across two runs (same rPr/lang=\"pl-PL\"):\n#   run 1: \"...sekwencji\"\n#   run 2: \" oraz chromosomie.\"   (xml:space=\"preserve\")\n\n$d = $word.ActiveDocument\n\n$oldText = \"bazuj\u0105c na pocz\u0105tku i ko\u0144cu sekwencji.\"\n$newText = \"bazuj\u0105c na pocz\u0105tku i ko\u0144cu sekwencji oraz chromosomie.\"\n\n# Step 1: replace the old sentence ending with the extended one in a single\n# Find/Replace call so the run keeps uniform formatting (this naturally\n# collapses into one merged run, same as typing a replacement would).\n$find = $d.Content.Find\n$found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\nif (-not $found) {\n    throw \"Could not find the target sentence to replace.\"\n}\n\n# Step 2: re-select just the newly appended suffix (\" oraz chromosomie.\")\n# and nudge a character-formatting property on/off. This forces the host to\n# materialize it as its own run (matching the target OOXML: two runs with\n# identical rPr) instead of silently leaving it merged with the previous run.\n$find2 = $d.Content.Find\n$found2 = $find2.Execute(\" oraz chromosomie.\")\nif (-not $found2) {\n    throw \"Could not find the newly inserted suffix.\"\n}\n$suffixRange = $find2.Parent\n$suffixRange.Font.Bold = 1\n$suffixRange.Font.Bold = 0\n"}
